# Add three product/colour-choice labels to the "Product Sheet" worksheet.
# (Per the commit message: "Made a few colour choices for the product
# items" — these are placeholder labels for the colour profiles that will
# later be wired up to the colour-picking screen.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "FUCK"
$ws.Range("A2").Value = "YEAH"
$ws.Range("A3").Value = "NO"

# Typing the values and pressing Enter in real Excel leaves the selection
# on the cell below the last entry - move the active selection to A4 to
# match.
$ws.Range("A4").Select()
